$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.909.76"
$ws.Range("E2").Value = "'  -1.38%  "
$ws.Range("D3").Value = "'2.223.92"
$ws.Range("E3").Value = "'  -1.13%  "
$ws.Range("E4").Value = "'  -0.30%  "
$ws.Range("D5").Value = "'298.19"
$ws.Range("E5").Value = "'  -3.50%  "
$ws.Range("D6").Value = "'89.78"
$ws.Range("E6").Value = "'  -5.50%  "
$ws.Range("D7").Value = "'0.556"
$ws.Range("E7").Value = "'  -3.05%  "
$ws.Range("E8").Value = "'  -0.32%  "
$ws.Range("D9").Value = "'0.489"
$ws.Range("E9").Value = "'  -7.41%  "
$ws.Range("D10").Value = "'32.72"
$ws.Range("E10").Value = "'  -6.80%  "
$ws.Range("D11").Value = "'0.0775"
$ws.Range("E11").Value = "'  -4.18%  "
$ws.Range("D12").Value = "'6.93"
$ws.Range("E12").Value = "'  -5.01%  "
$ws.Range("D13").Value = "'0.103"
$ws.Range("E13").Value = "'  -0.92%  "
$ws.Range("D14").Value = "'2.562.37"
$ws.Range("E14").Value = "'  -1.20%  "
$ws.Range("D15").Value = "'2.217.11"
$ws.Range("E15").Value = "'  -4.06%  "
$ws.Range("D16").Value = "'13.46"
$ws.Range("E16").Value = "'  -1.49%  "
$ws.Range("D17").Value = "'0.774"
$ws.Range("E17").Value = "'  -7.74%  "
$ws.Range("D18").Value = "'43.975.84"
$ws.Range("E18").Value = "'  -0.55%  "
$ws.Range("D19").Value = "'0.0₃0900"
$ws.Range("E19").Value = "'  -6.84%  "
$ws.Range("D20").Value = "'5.90"
$ws.Range("E20").Value = "'  -7.88%  "
$ws.Range("D21").Value = "'10.92"
$ws.Range("E21").Value = "'  -10.15%  "
$ws.Range("D22").Value = "'64.63"
$ws.Range("E22").Value = "'  -1.98%  "
$ws.Range("D23").Value = "'235.31"
$ws.Range("E23").Value = "'  -1.30%  "
$ws.Range("E24").Value = "'  -7.18%  "
$ws.Range("E25").Value = "'  +0.52%  "
$ws.Range("D26").Value = "'1.84"
$ws.Range("E26").Value = "'  -8.20%  "
$ws.Range("E27").Value = "'  +0.51%  "
$ws.Range("D28").Value = "'37.80"
$ws.Range("E28").Value = "'  +0.04%  "
$ws.Range("D29").Value = "'9.29"
$ws.Range("E29").Value = "'  -5.82%  "
$ws.Range("D30").Value = "'19.24"
$ws.Range("E30").Value = "'  -4.32%  "
$ws.Range("D31").Value = "'148.36"
$ws.Range("E31").Value = "'  -2.74%  "
$ws.Range("D32").Value = "'5.38"
$ws.Range("E32").Value = "'  -10.44%  "
$ws.Range("B33").Value = "'Hedera"
$ws.Range("C33").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0749"
$ws.Range("E33").Value = "'  -6.67%  "
$ws.Range("B34").Value = "'WEMIXToken"
$ws.Range("C34").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.51"
$ws.Range("E34").Value = "'  -4.71%  "
$ws.Range("E35").Value = "'  -4.31%  "
$ws.Range("B36").Value = "'LidoDAOToken"
$ws.Range("C36").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.79"
$ws.Range("E36").Value = "'  -12.12%  "
$ws.Range("B37").Value = "'Kaspa"
$ws.Range("C37").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.101"
$ws.Range("E37").Value = "'  -7.27%  "
$ws.Range("E38").Value = "'  -7.24%  "
$ws.Range("D39").Value = "'0.0300"
$ws.Range("E39").Value = "'  -1.13%  "
$ws.Range("E40").Value = "'  -8.25%  "
$ws.Range("E41").Value = "'  -8.35%  "
$ws.Range("E42").Value = "'  -0.36%  "
$ws.Range("D43").Value = "'13.03"
$ws.Range("E43").Value = "'  -10.09%  "
$ws.Range("D44").Value = "'1.802.96"
$ws.Range("E44").Value = "'  +2.77%  "
$ws.Range("D45").Value = "'1.76"
$ws.Range("E45").Value = "'  +10.43%  "
$ws.Range("D46").Value = "'0.177"
$ws.Range("E46").Value = "'  -8.29%  "
$ws.Range("D47").Value = "'73.38"
$ws.Range("E47").Value = "'  -9.32%  "
$ws.Range("B48").Value = "'EnergySwap"
$ws.Range("C48").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'14.02"
$ws.Range("E48").Value = "'  -2.56%  "
$ws.Range("B49").Value = "'Aave"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'93.70"
$ws.Range("E49").Value = "'  -6.33%  "
$ws.Range("D50").Value = "'66.43"
$ws.Range("E50").Value = "'  -6.51%  "
$ws.Range("D51").Value = "'2.443.62"
$ws.Range("E51").Value = "'  -1.20%  "
